$d = $word.ActiveDocument

$replacements = @(
    @("2025-09-25 Thursday", "2025-09-26 Friday"),
    @("441×4=", "742×3="),
    @("488×7=", "246×4="),
    @("191×9=", "844×3="),
    @("288×4=", "724×8="),
    @("186×3=", "147×3="),
    @("972×6=", "567×9="),
    @("666×5=", "879×3="),
    @("716×4=", "999×8="),
    @("199×9=", "742×5="),
    @("535×5=", "571×3="),
    @("828×8=", "312×5="),
    @("888×3=", "457×3="),
    @("384×8=", "651×9="),
    @("280×6=", "301×4="),
    @("482×4=", "307×6="),
    @("795×3=", "658×4="),
    @("883×6=", "692×7="),
    @("903×5=", "451×7="),
    @("304×2=", "619×5="),
    @("427×5=", "187×6="),
    @("397×4=", "418×9="),
    @("287×6=", "237×4="),
    @("938×4=", "940×6="),
    @("162×2=", "968×6="),
    @("464×3=", "140×7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
